$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.318.50"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -3.90%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.983.54"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.45%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "535.24"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.62%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.69"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.95%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.977.17"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.44%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.495"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.58%  "

$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.148"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -5.95%  "

$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.11"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.98%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.65%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.04%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.467.27"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.36%  "

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.03%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.357.26"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.81%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.985.36"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.38%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.19%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "464.32"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -5.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.13"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.667"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -5.42%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.73%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.20"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.87"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.25%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.16%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.53%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.67"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -7.56%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.03%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.40%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.13%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.44"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.31%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "54.94"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.94%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -6.39%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.84"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.92%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "454.20"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -9.69%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.143.26"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.45%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.43%  "

$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.119"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.49%  "

$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0381"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.31%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.76%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.42"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -9.99%  "

$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.09%  "

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "26.11"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.71%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -7.30%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -6.38%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.37%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "117.31"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.10%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -9.58%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +6.36%  "
